$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.Value = "'29.242.15"
$d.Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "
$d = $ws.Range("D3")
$d.Value = "'1.842.94"
$d.Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.03%  "
$d = $ws.Range("D5")
$d.Value = "'240.44"
$d.Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$d = $ws.Range("D6")
$d.Value = "'0.6744"
$d.Style = "Normal"
$ws.Range("E6").Value = "  -1.75%  "
$d = $ws.Range("D7")
$d.Value = "'0.9998"
$d.Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$d = $ws.Range("D8")
$d.Value = "'0.07427"
$d.Style = "Normal"
$ws.Range("E8").Value = "  -0.72%  "
$d = $ws.Range("D9")
$d.Value = "'0.2949"
$d.Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "
$d = $ws.Range("D10")
$d.Value = "'22.83"
$d.Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$d = $ws.Range("D11")
$d.Value = "'0.07717"
$d.Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$d = $ws.Range("D12")
$d.Value = "'1.857.25"
$d.Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "
$d = $ws.Range("D13")
$d.Value = "'5.001"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -1.16%  "
$d = $ws.Range("D14")
$d.Value = "'0.6703"
$d.Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "
$d = $ws.Range("D15")
$d.Value = "'86.06"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -1.86%  "
$d = $ws.Range("D16")
$d.Value = "'6.122"
$d.Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$d = $ws.Range("D17")
$d.Value = "'29.286.00"
$d.Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$d = $ws.Range("D18")
$d.Value = "'0.000008312"
$d.Style = "Normal"
$ws.Range("E18").Value = "  +1.64%  "
$d = $ws.Range("D19")
$d.Value = "'228.53"
$d.Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.14%  "
$d = $ws.Range("D22")
$d.Value = "'7.189"
$d.Style = "Normal"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("E23").Value = "  +0.10%  "
$d = $ws.Range("D24")
$d.Value = "'160.72"
$d.Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "
$d = $ws.Range("D25")
$d.Value = "'8.695"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$d = $ws.Range("D26")
$d.Value = "'0.1403"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -3.72%  "
$d = $ws.Range("D27")
$d.Value = "'17.99"
$d.Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$d = $ws.Range("D28")
$d.Value = "'1.507"
$d.Style = "Normal"
$ws.Range("E28").Value = "  -0.35%  "
$d = $ws.Range("D29")
$d.Value = "'4.176"
$d.Style = "Normal"
$ws.Range("E29").Value = "  -2.39%  "
$d = $ws.Range("D30")
$d.Value = "'4.069"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$d = $ws.Range("D31")
$d.Value = "'1.190"
$d.Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$d = $ws.Range("D32")
$d.Value = "'0.05310"
$d.Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "
$d = $ws.Range("D33")
$d.Value = "'0.7614"
$d.Style = "Normal"
$ws.Range("E33").Value = "  -0.56%  "
$d = $ws.Range("D34")
$d.Value = "'1.871"
$d.Style = "Normal"
$ws.Range("E34").Value = "  +1.35%  "
$d = $ws.Range("D35")
$d.Value = "'1.135"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$d = $ws.Range("D36")
$d.Value = "'2.676"
$d.Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "
$d = $ws.Range("D37")
$d.Value = "'1.329.08"
$d.Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$d = $ws.Range("D38")
$d.Value = "'0.01805"
$d.Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "
$d = $ws.Range("D39")
$d.Value = "'2.723"
$d.Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$d = $ws.Range("D40")
$d.Value = "'0.9213"
$d.Style = "Normal"
$ws.Range("E40").Value = "  -1.31%  "
$d = $ws.Range("D41")
$d.Value = "'5.934"
$d.Style = "Normal"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  +0.22%  "
$d = $ws.Range("D43")
$d.Value = "'103.48"
$d.Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "
$d = $ws.Range("D44")
$d.Value = "'0.08134"
$d.Style = "Normal"
$ws.Range("E44").Value = "  +14.71%  "
$d = $ws.Range("D45")
$d.Value = "'0.00000000127"
$d.Style = "Normal"
$ws.Range("E45").Value = "  +3.00%  "
$d = $ws.Range("D46")
$d.Value = "'1.999.52"
$d.Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "
$d = $ws.Range("D47")
$d.Value = "'0.5160"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
$d = $ws.Range("D48")
$d.Value = "'1.778"
$d.Style = "Normal"
$ws.Range("E48").Value = "  +0.31%  "
$d = $ws.Range("D49")
$d.Value = "'63.76"
$d.Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "
$d = $ws.Range("D50")
$d.Value = "'9.136"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -4.32%  "
$d = $ws.Range("D51")
$d.Value = "'0.05956"
$d.Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
